{"js": "// Update the worksheet date and every \"A\u00d7B=C\" answer cell to the new\n// values from the commit. Every <w:t> run in the document changes, and\n// every old value is unique, so a straightforward search-and-replace\n// (old text -> new text) for each pair reproduces the diff exactly.\nconst replacements = [\n  [\"2024-02-18 Sunday\", \"2024-02-19 Monday\"],\n  [\"319\u00d73=957\", \"296\u00d74=1184\"],\n  [\"263\u00d79=2367\", \"235\u00d77=1645\"],\n  [\"761\u00d77=5327\", \"940\u00d79=8460\"],\n  [\"681\u00d77=4767\", \"205\u00d73=615\"],\n  [\"674\u00d76=4044\", \"588\u00d74=2352\"],\n  [\"961\u00d79=8649\", \"292\u00d74=1168\"],\n  [\"188\u00d77=1316\", \"518\u00d79=4662\"],\n  [\"816\u00d78=6528\", \"318\u00d74=1272\"],\n  [\"435\u00d76=2610\", \"131\u00d75=655\"],\n  [\"958\u00d75=4790\", \"746\u00d78=5968\"],\n  [\"987\u00d76=5922\", \"464\u00d72=928\"],\n  [\"503\u00d74=2012\", \"585\u00d78=4680\"],\n  [\"290\u00d78=2320\", \"193\u00d75=965\"],\n  [\"257\u00d79=2313\", \"795\u00d76=4770\"],\n  [\"492\u00d78=3936\", \"873\u00d74=3492\"],\n  [\"526\u00d74=2104\", \"386\u00d75=1930\"],\n  [\"127\u00d76=762\", \"438\u00d79=3942\"],\n  [\"396\u00d72=792\", \"828\u00d74=3312\"],\n  [\"149\u00d75=745\", \"489\u00d77=3423\"],\n  [\"922\u00d72=1844\", \"759\u00d74=3036\"],\n  [\"697\u00d76=4182\", \"384\u00d79=3456\"],\n  [\"151\u00d76=906\", \"153\u00d78=1224\"],\n  [\"356\u00d74=1424\", \"507\u00d77=3549\"],\n  [\"416\u00d77=2912\", \"505\u00d77=3535\"],\n  [\"797\u00d74=3188\", \"464\u00d75=2320\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every \"A\u00d7B=C\" answer cell to the new\n# values from the commit. Every run of text in the document changes, and\n# every old value is unique, so a plain Find/Replace (MatchWholeWord off,\n# wildcards off, Replace:=wdReplaceAll) per pair reproduces the diff\n# exactly.\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # Forward, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Wrap:=wdFindContinue(1), Format, ReplaceWith,\n    # Replace:=wdReplaceAll(2)\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-All '2024-02-18 Sunday' '2024-02-19 Monday'\nReplace-All '319\u00d73=957' '296\u00d74=1184'\nReplace-All '263\u00d79=2367' '235\u00d77=1645'\nReplace-All '761\u00d77=5327' '940\u00d79=8460'\nReplace-All '681\u00d77=4767' '205\u00d73=615'\nReplace-All '674\u00d76=4044' '588\u00d74=2352'\nReplace-All '961\u00d79=8649' '292\u00d74=1168'\nReplace-All '188\u00d77=1316' '518\u00d79=4662'\nReplace-All '816\u00d78=6528' '318\u00d74=1272'\nReplace-All '435\u00d76=2610' '131\u00d75=655'\nReplace-All '958\u00d75=4790' '746\u00d78=5968'\nReplace-All '987\u00d76=5922' '464\u00d72=928'\nReplace-All '503\u00d74=2012' '585\u00d78=4680'\nReplace-All '290\u00d78=2320' '193\u00d75=965'\nReplace-All '257\u00d79=2313' '795\u00d76=4770'\nReplace-All '492\u00d78=3936' '873\u00d74=3492'\nReplace-All '526\u00d74=2104' '386\u00d75=1930'\nReplace-All '127\u00d76=762' '438\u00d79=3942'\nReplace-All '396\u00d72=792' '828\u00d74=3312'\nReplace-All '149\u00d75=745' '489\u00d77=3423'\nReplace-All '922\u00d72=1844' '759\u00d74=3036'\nReplace-All '697\u00d76=4182' '384\u00d79=3456'\nReplace-All '151\u00d76=906' '153\u00d78=1224'\nReplace-All '356\u00d74=1424' '507\u00d77=3549'\nReplace-All '416\u00d77=2912' '505\u00d77=3535'\nReplace-All '797\u00d74=3188' '464\u00d75=2320'\n"}
